# V 0.50-B47 - SU9 checking - Add Control Surfaces Item as Graphics - SU9 Compatibility Check - Update QuickGuide
#
# Tabelle2 (sheet2) gains a new column at EB: a "SURF_ANI" marker column
# (mirrors the existing END_OF_COL / "X" marker pattern), which pushes the
# old EB ("END_OF_COL" header / "X" data marker) and EC (Title header /
# Tabelle1!AS.. formula) columns one place to the right, to EC and ED.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert a brand-new column before EB; this shifts the old EB -> EC and
# the old EC -> ED, carrying their formatting/formulas along with them.
$ws.Columns("EB").Insert()

# Header row: new EB1 marker string (added to sharedStrings.xml).
$ws.Range("EB1").Value = "SURF_ANI"

# Data rows 2-40: new EB column uses the same "|" marker used throughout
# the rest of the sheet's data columns.
$ws.Range("EB2:EB40").Value = "|"

# Leave the selection where the author ended up after the edit.
$ws.Range("EH6").Select()
